$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu values (380 kV case) for rows 2-25, columns B:F and I:N
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.049461418104183
$ws.Cells.Item(2, 4).Value = 1.054268500874967
$ws.Cells.Item(2, 5).Value = 1.056423404658036
$ws.Cells.Item(2, 6).Value = 1.066577067320107
$ws.Cells.Item(2, 9).Value = 1.044679309074269
$ws.Cells.Item(2, 10).Value = 1.054499741971103
$ws.Cells.Item(2, 11).Value = 1.057012178741535
$ws.Cells.Item(2, 12).Value = 1.059161156060306
$ws.Cells.Item(2, 13).Value = 1.069287242618522
$ws.Cells.Item(2, 14).Value = 1.021966965282052
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.050416860626251
$ws.Cells.Item(3, 4).Value = 1.055010325809476
$ws.Cells.Item(3, 5).Value = 1.05726754969049
$ws.Cells.Item(3, 6).Value = 1.067488695386129
$ws.Cells.Item(3, 9).Value = 1.044916980563474
$ws.Cells.Item(3, 10).Value = 1.055104444839172
$ws.Cells.Item(3, 11).Value = 1.057567395799492
$ws.Cells.Item(3, 12).Value = 1.059818859499298
$ws.Cells.Item(3, 13).Value = 1.07001425038649
$ws.Cells.Item(3, 14).Value = 1.02217139681673
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.051035639155037
$ws.Cells.Item(4, 4).Value = 1.05549081568895
$ws.Cells.Item(4, 5).Value = 1.057814618097268
$ws.Cells.Item(4, 6).Value = 1.068079551094568
$ws.Cells.Item(4, 9).Value = 1.045069853785605
$ws.Cells.Item(4, 10).Value = 1.055495632441421
$ws.Cells.Item(4, 11).Value = 1.057926469216746
$ws.Cells.Item(4, 12).Value = 1.060244639408543
$ws.Cells.Item(4, 13).Value = 1.070485008083912
$ws.Cells.Item(4, 14).Value = 1.022303563609145
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.051295902651391
$ws.Cells.Item(5, 4).Value = 1.05569292702714
$ws.Cells.Item(5, 5).Value = 1.058044807743626
$ws.Cells.Item(5, 6).Value = 1.068328177888631
$ws.Cells.Item(5, 9).Value = 1.045133901802177
$ws.Cells.Item(5, 10).Value = 1.055660063799455
$ws.Cells.Item(5, 11).Value = 1.058077377358864
$ws.Cells.Item(5, 12).Value = 1.060423684454691
$ws.Cells.Item(5, 13).Value = 1.070682993573444
$ws.Cells.Item(5, 14).Value = 1.022359098794294
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.05133960957142
$ws.Cells.Item(6, 4).Value = 1.055726869018573
$ws.Cells.Item(6, 5).Value = 1.058083469407039
$ws.Cells.Item(6, 6).Value = 1.068369936926969
$ws.Cells.Item(6, 9).Value = 1.045144642835417
$ws.Cells.Item(6, 10).Value = 1.055687671116139
$ws.Cells.Item(6, 11).Value = 1.058102712753365
$ws.Cells.Item(6, 12).Value = 1.060453749638874
$ws.Cells.Item(6, 13).Value = 1.070716240801598
$ws.Cells.Item(6, 14).Value = 1.02236842175759
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.051039116303048
$ws.Cells.Item(7, 4).Value = 1.055493515866755
$ws.Cells.Item(7, 5).Value = 1.057817693110315
$ws.Cells.Item(7, 6).Value = 1.068082872351948
$ws.Cells.Item(7, 9).Value = 1.045070710463038
$ws.Cells.Item(7, 10).Value = 1.055497829676504
$ws.Cells.Item(7, 11).Value = 1.057928485842043
$ws.Cells.Item(7, 12).Value = 1.060247031633283
$ws.Cells.Item(7, 13).Value = 1.070487653267055
$ws.Cells.Item(7, 14).Value = 1.022304305782513
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.049784201232633
$ws.Cells.Item(8, 4).Value = 1.054519103944159
$ws.Cells.Item(8, 5).Value = 1.056708510538015
$ws.Cells.Item(8, 6).Value = 1.066884954258747
$ws.Cells.Item(8, 9).Value = 1.044759820616519
$ws.Cells.Item(8, 10).Value = 1.054704123311478
$ws.Cells.Item(8, 11).Value = 1.057199855483457
$ws.Cells.Item(8, 12).Value = 1.059383387545857
$ws.Cells.Item(8, 13).Value = 1.069532868163518
$ws.Cells.Item(8, 14).Value = 1.022036077242918
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.047577089175586
$ws.Cells.Item(9, 4).Value = 1.052805800683007
$ws.Cells.Item(9, 5).Value = 1.05476056341454
$ws.Cells.Item(9, 6).Value = 1.06478157547254
$ws.Cells.Item(9, 9).Value = 1.044205001552015
$ws.Cells.Item(9, 10).Value = 1.053304824539093
$ws.Cells.Item(9, 11).Value = 1.05591451306999
$ws.Cells.Item(9, 12).Value = 1.057863132298486
$ws.Cells.Item(9, 13).Value = 1.067853036663542
$ws.Cells.Item(9, 14).Value = 1.021562569258168
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.046108573588409
$ws.Cells.Item(10, 4).Value = 1.051666193334444
$ws.Cells.Item(10, 5).Value = 1.053466437071522
$ws.Cells.Item(10, 6).Value = 1.063384458296338
$ws.Cells.Item(10, 9).Value = 1.043830454947992
$ws.Cells.Item(10, 10).Value = 1.052371557544721
$ws.Cells.Item(10, 11).Value = 1.055056736644592
$ws.Cells.Item(10, 12).Value = 1.056850769285814
$ws.Cells.Item(10, 13).Value = 1.066734986353366
$ws.Cells.Item(10, 14).Value = 1.021246347423532
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.045473389060422
$ws.Cells.Item(11, 4).Value = 1.05117336372091
$ws.Cells.Item(11, 5).Value = 1.052907152362461
$ws.Cells.Item(11, 6).Value = 1.06278072771697
$ws.Cells.Item(11, 9).Value = 1.043667172073276
$ws.Cells.Item(11, 10).Value = 1.051967360075136
$ws.Cells.Item(11, 11).Value = 1.054685113733517
$ws.Cells.Item(11, 12).Value = 1.056412689487506
$ws.Cells.Item(11, 13).Value = 1.066251309579878
$ws.Cells.Item(11, 14).Value = 1.021109294681079
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.045237558263516
$ws.Cells.Item(12, 4).Value = 1.050990400586583
$ws.Cells.Item(12, 5).Value = 1.052699572872225
$ws.Cells.Item(12, 6).Value = 1.062556661627029
$ws.Cells.Item(12, 9).Value = 1.043606356555163
$ws.Cells.Item(12, 10).Value = 1.051817211097766
$ws.Cells.Item(12, 11).Value = 1.054547047389661
$ws.Cells.Item(12, 12).Value = 1.056250010331657
$ws.Cells.Item(12, 13).Value = 1.06607171878703
$ws.Cells.Item(12, 14).Value = 1.021058368580856
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.045288140001717
$ws.Cells.Item(13, 4).Value = 1.051029642445725
$ws.Cells.Item(13, 5).Value = 1.052744091944643
$ws.Cells.Item(13, 6).Value = 1.062604716108433
$ws.Cells.Item(13, 9).Value = 1.043619409147719
$ws.Cells.Item(13, 10).Value = 1.051849419092099
$ws.Cells.Item(13, 11).Value = 1.054576664387376
$ws.Cells.Item(13, 12).Value = 1.056284903590684
$ws.Cells.Item(13, 13).Value = 1.066110238511375
$ws.Cells.Item(13, 14).Value = 1.02106929323774
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.045453893065082
$ws.Cells.Item(14, 4).Value = 1.051158237970276
$ws.Cells.Item(14, 5).Value = 1.05288999042835
$ws.Cells.Item(14, 6).Value = 1.062762202534481
$ws.Cells.Item(14, 9).Value = 1.043662148406193
$ws.Cells.Item(14, 10).Value = 1.051954948953734
$ws.Cells.Item(14, 11).Value = 1.054673701713198
$ws.Cells.Item(14, 12).Value = 1.056399241490313
$ws.Cells.Item(14, 13).Value = 1.066236463148558
$ws.Cells.Item(14, 14).Value = 1.021105085490997
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.04555603294368
$ws.Cells.Item(15, 4).Value = 1.051237482643659
$ws.Cells.Item(15, 5).Value = 1.05297990504059
$ws.Cells.Item(15, 6).Value = 1.062859259853834
$ws.Cells.Item(15, 9).Value = 1.043688459627366
$ws.Cells.Item(15, 10).Value = 1.052019967818492
$ws.Cells.Item(15, 11).Value = 1.054733485793305
$ws.Cells.Item(15, 12).Value = 1.056469694599796
$ws.Cells.Item(15, 13).Value = 1.066314243400296
$ws.Cells.Item(15, 14).Value = 1.021127135825894
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.046150743366402
$ws.Cells.Item(16, 4).Value = 1.05169891416069
$ws.Cells.Item(16, 5).Value = 1.053503577842082
$ws.Cells.Item(16, 6).Value = 1.063424551955312
$ws.Cells.Item(16, 9).Value = 1.04384126832803
$ws.Cells.Item(16, 10).Value = 1.052398381056952
$ws.Cells.Item(16, 11).Value = 1.055081395904314
$ws.Cells.Item(16, 12).Value = 1.056879849194626
$ws.Cells.Item(16, 13).Value = 1.066767095897858
$ws.Cells.Item(16, 14).Value = 1.021255440530411
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.04652397557882
$ws.Cells.Item(17, 4).Value = 1.051988527048198
$ws.Cells.Item(17, 5).Value = 1.053832354243824
$ws.Cells.Item(17, 6).Value = 1.063779475205599
$ws.Cells.Item(17, 9).Value = 1.043936826601684
$ws.Cells.Item(17, 10).Value = 1.052635727166339
$ws.Cells.Item(17, 11).Value = 1.055299577875203
$ws.Cells.Item(17, 12).Value = 1.0571372041006
$ws.Cells.Item(17, 13).Value = 1.067051278685228
$ws.Cells.Item(17, 14).Value = 1.02133588907711
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.046741742354647
$ws.Cells.Item(18, 4).Value = 1.05215751383425
$ws.Cells.Item(18, 5).Value = 1.054024228117785
$ws.Cells.Item(18, 6).Value = 1.063986614547038
$ws.Cells.Item(18, 9).Value = 1.043992457789727
$ws.Cells.Item(18, 10).Value = 1.052774158704597
$ws.Cells.Item(18, 11).Value = 1.05542682035701
$ws.Cells.Item(18, 12).Value = 1.057287341811005
$ws.Cells.Item(18, 13).Value = 1.067217080684454
$ws.Cells.Item(18, 14).Value = 1.021382801110553
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.046816006487952
$ws.Cells.Item(19, 4).Value = 1.052215144178795
$ws.Cells.Item(19, 5).Value = 1.054089669774888
$ws.Cells.Item(19, 6).Value = 1.064057263783622
$ws.Cells.Item(19, 9).Value = 1.044011408534556
$ws.Cells.Item(19, 10).Value = 1.052821358821623
$ws.Cells.Item(19, 11).Value = 1.05547020343885
$ws.Cells.Item(19, 12).Value = 1.05733853942743
$ws.Cells.Item(19, 13).Value = 1.06727362213544
$ws.Cells.Item(19, 14).Value = 1.021398794812555
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.046483924384149
$ws.Cells.Item(20, 4).Value = 1.051957448059536
$ws.Cells.Item(20, 5).Value = 1.053797068863312
$ws.Cells.Item(20, 6).Value = 1.063741383020277
$ws.Cells.Item(20, 9).Value = 1.043926585098279
$ws.Cells.Item(20, 10).Value = 1.052610263040028
$ws.Cells.Item(20, 11).Value = 1.055276171010866
$ws.Cells.Item(20, 12).Value = 1.057109589559402
$ws.Cells.Item(20, 13).Value = 1.067020784106333
$ws.Cells.Item(20, 14).Value = 1.021327258972761
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.045405080012193
$ws.Cells.Item(21, 4).Value = 1.051120367139711
$ws.Cells.Item(21, 5).Value = 1.052847022429576
$ws.Cells.Item(21, 6).Value = 1.062715821551801
$ws.Cells.Item(21, 9).Value = 1.043649567310008
$ws.Cells.Item(21, 10).Value = 1.051923873364186
$ws.Cells.Item(21, 11).Value = 1.054645127433098
$ws.Cells.Item(21, 12).Value = 1.05636557062758
$ws.Cells.Item(21, 13).Value = 1.066199291243753
$ws.Cells.Item(21, 14).Value = 1.021094546074919
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.044727375835591
$ws.Cells.Item(22, 4).Value = 1.050594615681381
$ws.Cells.Item(22, 5).Value = 1.05225063865407
$ws.Cells.Item(22, 6).Value = 1.062072089123504
$ws.Cells.Item(22, 9).Value = 1.043474440926496
$ws.Cells.Item(22, 10).Value = 1.051492243760829
$ws.Cells.Item(22, 11).Value = 1.054248197726587
$ws.Cells.Item(22, 12).Value = 1.055898026811126
$ws.Cells.Item(22, 13).Value = 1.065683181829957
$ws.Cells.Item(22, 14).Value = 1.020948122775075
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.045086581668773
$ws.Cells.Item(23, 4).Value = 1.050873273447375
$ws.Cells.Item(23, 5).Value = 1.052566702633567
$ws.Cells.Item(23, 6).Value = 1.062413241114482
$ws.Cells.Item(23, 9).Value = 1.043567369014204
$ws.Cells.Item(23, 10).Value = 1.051721065007085
$ws.Cells.Item(23, 11).Value = 1.054458633219989
$ws.Cells.Item(23, 12).Value = 1.056145856526806
$ws.Cells.Item(23, 13).Value = 1.065956743270638
$ws.Cells.Item(23, 14).Value = 1.021025754585842
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.046502021588888
$ws.Cells.Item(24, 4).Value = 1.051971491129986
$ws.Cells.Item(24, 5).Value = 1.053813012486235
$ws.Cells.Item(24, 6).Value = 1.063758594872825
$ws.Cells.Item(24, 9).Value = 1.043931213121437
$ws.Cells.Item(24, 10).Value = 1.052621769208762
$ws.Cells.Item(24, 11).Value = 1.055286747625274
$ws.Cells.Item(24, 12).Value = 1.057122067298763
$ws.Cells.Item(24, 13).Value = 1.067034563161214
$ws.Cells.Item(24, 14).Value = 1.021331158583308
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.048147175351508
$ws.Cells.Item(25, 4).Value = 1.053248279409795
$ws.Cells.Item(25, 5).Value = 1.055263366608111
$ws.Cells.Item(25, 6).Value = 1.065324450672385
$ws.Cells.Item(25, 9).Value = 1.044349260666548
$ws.Cells.Item(25, 10).Value = 1.053666651961857
$ws.Cells.Item(25, 11).Value = 1.056246964775026
$ws.Cells.Item(25, 12).Value = 1.058255958645845
$ws.Cells.Item(25, 13).Value = 1.068286994992095
$ws.Cells.Item(25, 14).Value = 1.021685081169173
